$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: "_old" columns become "_FV2404", "_new" columns become "_FV2410"
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Create a table (ListObject) over the data range, with headers
$range = $ws.Range("A1:U78")
$listObject = $ws.ListObjects.Add(1, $range, 0, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""
$listObject.ShowTableStyleRowStripes = $true
$listObject.ShowTableStyleColumnStripes = $false
$listObject.ShowTableStyleFirstColumn = $false
$listObject.ShowTableStyleLastColumn = $false

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
